$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 = 9, C3 = 80 (plain numbers)
$ws.Range("C2").Value = 9
$ws.Range("C3").Value = 80

# D3 = "true" stored as text (shared string), not as a Boolean.
# Typing a formula that evaluates to the text "true" and then pasting the
# result back as a value keeps Excel from auto-boxing it into a Boolean,
# and avoids leaving a stray quote-prefix style behind.
$ws.Range("D3").Formula = "=""true"""
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)  # xlPasteValues

# D4: empty cell, formatted with an underlined font (new style).
$ws.Range("D4").Font.Underline = $true

# Page setup: paper size 9 (A4), portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the final selection on D4, matching the recorded view state.
$null = $ws.Range("D4").Select()
